$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("BVT")
$ws1.PageSetup.PaperSize = 9
$ws1.PageSetup.Orientation = 1
$ws1.PageSetup.PrintQuality = 0
